$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string for row 11 Tasks done + row 12 new entry
$ws.Range("B11").Value = "improvments on the inventory and a beginning on some puzzle mechanics"
$ws.Range("C11").Value = 5

$ws.Range("A12").Value = 45628
$ws.Range("A12").NumberFormat = "d-mmm"

$ws.Range("A11").Select()
